{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" practice answers in the table\n// with the new values from the commit, matching each cell's exact old text\n// and writing back the new text (format-preserving, in place).\nconst replacements = [\n  [\"153\u00f76=25, 3\", \"204\u00f79=22, 6\"],\n  [\"482\u00f79=53, 5\", \"593\u00f79=65, 8\"],\n  [\"536\u00f77=76, 4\", \"606\u00f76=101, 0\"],\n  [\"659\u00f76=109, 5\", \"601\u00f73=200, 1\"],\n  [\"913\u00f74=228, 1\", \"345\u00f72=172, 1\"],\n  [\"842\u00f72=421, 0\", \"742\u00f74=185, 2\"],\n  [\"660\u00f76=110, 0\", \"564\u00f73=188, 0\"],\n  [\"836\u00f74=209, 0\", \"870\u00f75=174, 0\"],\n  [\"229\u00f77=32, 5\", \"596\u00f76=99, 2\"],\n  [\"122\u00f79=13, 5\", \"811\u00f72=405, 1\"],\n  [\"304\u00f73=101, 1\", \"540\u00f75=108, 0\"],\n  [\"631\u00f72=315, 1\", \"947\u00f79=105, 2\"],\n  [\"479\u00f72=239, 1\", \"483\u00f76=80, 3\"],\n  [\"428\u00f72=214, 0\", \"802\u00f74=200, 2\"],\n  [\"976\u00f78=122, 0\", \"371\u00f79=41, 2\"],\n  [\"748\u00f72=374, 0\", \"878\u00f73=292, 2\"],\n  [\"579\u00f75=115, 4\", \"177\u00f78=22, 1\"],\n  [\"838\u00f74=209, 2\", \"829\u00f73=276, 1\"],\n  [\"735\u00f79=81, 6\", \"589\u00f75=117, 4\"],\n  [\"713\u00f74=178, 1\", \"695\u00f74=173, 3\"],\n  [\"393\u00f77=56, 1\", \"185\u00f74=46, 1\"],\n  [\"281\u00f76=46, 5\", \"450\u00f74=112, 2\"],\n  [\"860\u00f76=143, 2\", \"454\u00f75=90, 4\"],\n  [\"407\u00f75=81, 2\", \"613\u00f75=122, 3\"],\n  [\"429\u00f77=61, 2\", \"873\u00f73=291, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" practice answers in the table\n# with the new values from the commit, matching each cell's exact old text\n# and writing back the new text (format-preserving, in place).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"153\u00f76=25, 3\", \"204\u00f79=22, 6\"),\n    @(\"482\u00f79=53, 5\", \"593\u00f79=65, 8\"),\n    @(\"536\u00f77=76, 4\", \"606\u00f76=101, 0\"),\n    @(\"659\u00f76=109, 5\", \"601\u00f73=200, 1\"),\n    @(\"913\u00f74=228, 1\", \"345\u00f72=172, 1\"),\n    @(\"842\u00f72=421, 0\", \"742\u00f74=185, 2\"),\n    @(\"660\u00f76=110, 0\", \"564\u00f73=188, 0\"),\n    @(\"836\u00f74=209, 0\", \"870\u00f75=174, 0\"),\n    @(\"229\u00f77=32, 5\", \"596\u00f76=99, 2\"),\n    @(\"122\u00f79=13, 5\", \"811\u00f72=405, 1\"),\n    @(\"304\u00f73=101, 1\", \"540\u00f75=108, 0\"),\n    @(\"631\u00f72=315, 1\", \"947\u00f79=105, 2\"),\n    @(\"479\u00f72=239, 1\", \"483\u00f76=80, 3\"),\n    @(\"428\u00f72=214, 0\", \"802\u00f74=200, 2\"),\n    @(\"976\u00f78=122, 0\", \"371\u00f79=41, 2\"),\n    @(\"748\u00f72=374, 0\", \"878\u00f73=292, 2\"),\n    @(\"579\u00f75=115, 4\", \"177\u00f78=22, 1\"),\n    @(\"838\u00f74=209, 2\", \"829\u00f73=276, 1\"),\n    @(\"735\u00f79=81, 6\", \"589\u00f75=117, 4\"),\n    @(\"713\u00f74=178, 1\", \"695\u00f74=173, 3\"),\n    @(\"393\u00f77=56, 1\", \"185\u00f74=46, 1\"),\n    @(\"281\u00f76=46, 5\", \"450\u00f74=112, 2\"),\n    @(\"860\u00f76=143, 2\", \"454\u00f75=90, 4\"),\n    @(\"407\u00f75=81, 2\", \"613\u00f75=122, 3\"),\n    @(\"429\u00f77=61, 2\", \"873\u00f73=291, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
